# "PROVA - versione finale"
#
# On the "Release Date (RD)" sheet the small lookup table (id, release_date,
# tassativita) grows from 2 data rows to 4 data rows:
#   - the row that used to be row 2 (id 254187) is kept, but shifted down to
#     row 3;
#   - row 2 is replaced with a new entry (id 253974);
#   - a brand new row (id 254547) is inserted as the new row 4;
#   - the former row 3 (id 254967) ends up, unchanged, as the new row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Release Date (RD)")

# Insert two blank rows above the current last data row (row 3, id 254967)
# so the table has room for the extra entries. Inserting rows this way
# shifts the existing data down and keeps the row/column formatting
# (including the date number format on column B) intact.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(4).Insert()

# Row 2: new values overwrite the old row-2 entry.
$ws.Cells.Item(2, 1).Value = 253974
$ws.Cells.Item(2, 2).Value = 45980.58333333334
$ws.Cells.Item(2, 3).Value = 0

# Row 3: the original row-2 entry, now shifted down.
$ws.Cells.Item(3, 1).Value = 254187
$ws.Cells.Item(3, 2).Value = 45975.58333333334
$ws.Cells.Item(3, 3).Value = 0

# Row 4: brand new entry.
$ws.Cells.Item(4, 1).Value = 254547
$ws.Cells.Item(4, 2).Value = 45980.58333333334
$ws.Cells.Item(4, 3).Value = 0

# Row 5 already holds the original final row (id 254967) thanks to the two
# inserts above; it is left as-is.
